$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 18377640
$ws.Range("I32").Value = 38461784
$ws.Range("J32").Value = 3872426.8
$ws.Range("K32").Value = 38461784
$ws.Range("L32").Value = 3872426.8
$ws.Range("M32").Value = -38461458
$ws.Range("N32").Value = -3873078.8
$ws.Range("H62").Value = 9376.781000000001
$ws.Range("I62").Value = 11431.381
$ws.Range("J62").Value = 5454.364
$ws.Range("K62").Value = 11431.381
$ws.Range("L62").Value = 5454.364
$ws.Range("M62").Value = -10807.381
$ws.Range("N62").Value = -6702.364
$ws.Range("H65").Value = 9376.781000000001
$ws.Range("I65").Value = 11431.381
$ws.Range("J65").Value = 5454.364
$ws.Range("K65").Value = 57156.905
$ws.Range("L65").Value = 27271.82
$ws.Range("M65").Value = -54036.905
$ws.Range("N65").Value = -33511.82
$ws.Range("H121").Value = 1770
$ws.Range("I121").Value = 800
$ws.Range("J121").Value = 1908.5714
$ws.Range("K121").Value = 2400
$ws.Range("L121").Value = 5725.7142
$ws.Range("M121").Value = -653
$ws.Range("N121").Value = -9219.7142
$ws.Range("H132").Value = 3488.745
$ws.Range("I132").Value = 891.1627999999999
$ws.Range("J132").Value = 17450.75
$ws.Range("K132").Value = 2673.4884
$ws.Range("L132").Value = 52352.25
$ws.Range("M132").Value = -143.4883999999997
$ws.Range("N132").Value = -57412.25
$ws.Range("H135").Value = 233.33333
$ws.Range("I135").Value = 233.33333
$ws.Range("K135").Value = 2099.99997
$ws.Range("M135").Value = 435.0000300000002
$ws.Range("H137").Value = 2019801
$ws.Range("I137").Value = 2088717.5
$ws.Range("J137").Value = 1909534.8
$ws.Range("K137").Value = 6266152.5
$ws.Range("L137").Value = 5728604.4
$ws.Range("M137").Value = -6263602.5
$ws.Range("N137").Value = -5733704.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 878.75
$ws.Range("I45").Value = 806
$ws.Range("J45").Value = 951.5
$ws.Range("K45").Value = 806
$ws.Range("L45").Value = 951.5
$ws.Range("M45").Value = -429
$ws.Range("N45").Value = -1705.5
$ws.Range("H105").Value = 38370
$ws.Range("J105").Value = 38370
$ws.Range("L105").Value = 38370
$ws.Range("N105").Value = -45358
$ws.Range("H132").Value = 22022.84
$ws.Range("I132").Value = 32486.125
$ws.Range("J132").Value = 3421.4443
$ws.Range("K132").Value = 97458.375
$ws.Range("L132").Value = 10264.3329
$ws.Range("M132").Value = -94928.375
$ws.Range("N132").Value = -15324.3329

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 203.17647
$ws.Range("I22").Value = 193.77777
$ws.Range("J22").Value = 213.75
$ws.Range("K22").Value = 193.77777
$ws.Range("L22").Value = 213.75
$ws.Range("M22").Value = -20.77777
$ws.Range("N22").Value = -559.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2034.4286
$ws.Range("I22").Value = 2612.2
$ws.Range("J22").Value = 590
$ws.Range("K22").Value = 2612.2
$ws.Range("L22").Value = 590
$ws.Range("M22").Value = -2262.2
$ws.Range("N22").Value = -1290
$ws.Range("H31").Value = 2529779.2
$ws.Range("I31").Value = 1338.5
$ws.Range("J31").Value = 3502256.2
$ws.Range("K31").Value = 1338.5
$ws.Range("L31").Value = 3502256.2
$ws.Range("M31").Value = -1043.5
$ws.Range("N31").Value = -3502846.2
$ws.Range("H34").Value = 2529779.2
$ws.Range("I34").Value = 1338.5
$ws.Range("J34").Value = 3502256.2
$ws.Range("K34").Value = 1338.5
$ws.Range("L34").Value = 3502256.2
$ws.Range("M34").Value = -1136.5
$ws.Range("N34").Value = -3502660.2
$ws.Range("H58").Value = 4425.225
$ws.Range("I58").Value = 5802.95
$ws.Range("K58").Value = 5802.95
$ws.Range("M58").Value = -5599.95
$ws.Range("H132").Value = 9617081
$ws.Range("I132").Value = 13514651
$ws.Range("J132").Value = 3076
$ws.Range("K132").Value = 40543953
$ws.Range("L132").Value = 9228
$ws.Range("M132").Value = -40541423
$ws.Range("N132").Value = -14288
$ws.Range("H134").Value = 21740786
$ws.Range("I134").Value = 31250976
$ws.Range("J134").Value = 3206.2856
$ws.Range("K134").Value = 93752928
$ws.Range("L134").Value = 9618.856800000001
$ws.Range("M134").Value = -93750393
$ws.Range("N134").Value = -14688.8568
$ws.Range("H136").Value = 4425.225
$ws.Range("I136").Value = 5802.95
$ws.Range("K136").Value = 17408.85
$ws.Range("M136").Value = -14858.85

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 730.76086
$ws.Range("J5").Value = 1219
$ws.Range("L5").Value = 3657
$ws.Range("N5").Value = -3881
$ws.Range("H68").Value = 870.58826
$ws.Range("I68").Value = 522.53845
$ws.Range("J68").Value = 1086.0476
$ws.Range("K68").Value = 1567.61535
$ws.Range("L68").Value = 3258.142800000001
$ws.Range("M68").Value = -756.61535
$ws.Range("N68").Value = -4880.142800000001
$ws.Range("H71").Value = 870.58826
$ws.Range("I71").Value = 522.53845
$ws.Range("J71").Value = 1086.0476
$ws.Range("K71").Value = 4702.84605
$ws.Range("L71").Value = 9774.428400000001
$ws.Range("M71").Value = -646.8460500000001
$ws.Range("N71").Value = -17886.4284
$ws.Range("H131").Value = 1160.2
$ws.Range("I131").Value = 954.46155
$ws.Range("J131").Value = 1217.1063
$ws.Range("K131").Value = 2863.38465
$ws.Range("L131").Value = 3651.3189
$ws.Range("M131").Value = 2176.61535
$ws.Range("N131").Value = -13731.3189
$ws.Range("H132").Value = 1366.6842
$ws.Range("I132").Value = 2301
$ws.Range("J132").Value = 935.46155
$ws.Range("K132").Value = 20709
$ws.Range("L132").Value = 8419.15395
$ws.Range("M132").Value = -18179
$ws.Range("N132").Value = -13479.15395
$ws.Range("H135").Value = 730.76086
$ws.Range("J135").Value = 1219
$ws.Range("L135").Value = 10971
$ws.Range("N135").Value = -16041
$ws.Range("H138").Value = 1894
$ws.Range("I138").Value = 993.3333
$ws.Range("J138").Value = 10000
$ws.Range("K138").Value = 2979.9999
$ws.Range("L138").Value = 30000
$ws.Range("M138").Value = 2160.0001
$ws.Range("N138").Value = -40280

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 39399.2
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 39399.2
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 39399.2
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -41615.2
$ws.Range("H126").Value = 2531.3333
$ws.Range("I126").Value = 2180.4
$ws.Range("J126").Value = 2850.3635
$ws.Range("K126").Value = 6541.200000000001
$ws.Range("L126").Value = 8551.0905
$ws.Range("M126").Value = -4071.200000000001
$ws.Range("N126").Value = -13491.0905
$ws.Range("H136").Value = 23333.334
$ws.Range("J136").Value = 23333.334
$ws.Range("L136").Value = 70000.00199999999
$ws.Range("N136").Value = -75100.00199999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 875
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 966.6667
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 966.6667
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -1556.6667
$ws.Range("H27").Value = 875
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 966.6667
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 966.6667
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -1180.6667
$ws.Range("H40").Value = 1326.9231
$ws.Range("I40").Value = 1133.3334
$ws.Range("J40").Value = 1762.5
$ws.Range("K40").Value = 1133.3334
$ws.Range("L40").Value = 1762.5
$ws.Range("M40").Value = -997.3334
$ws.Range("N40").Value = -2034.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3493.5334
$ws.Range("I81").Value = 1283.5
$ws.Range("J81").Value = 4966.8887
$ws.Range("K81").Value = 2567
$ws.Range("L81").Value = 9933.777400000001
$ws.Range("M81").Value = -1506
$ws.Range("N81").Value = -12055.7774
$ws.Range("H84").Value = 3493.5334
$ws.Range("I84").Value = 1283.5
$ws.Range("J84").Value = 4966.8887
$ws.Range("K84").Value = 12835
$ws.Range("L84").Value = 49668.887
$ws.Range("M84").Value = -7531
$ws.Range("N84").Value = -60276.887
$ws.Range("H132").Value = 19609828
$ws.Range("I132").Value = 27028408
$ws.Range("J132").Value = 3576.2856
$ws.Range("K132").Value = 81085224
$ws.Range("L132").Value = 10728.8568
$ws.Range("M132").Value = -81082694
$ws.Range("N132").Value = -15788.8568
$ws.Range("H136").Value = 15165003
$ws.Range("I136").Value = 16309181
$ws.Range("J136").Value = 9253417
$ws.Range("K136").Value = 48927543
$ws.Range("L136").Value = 27760251
$ws.Range("M136").Value = -48924993
$ws.Range("N136").Value = -27765351

$wb.Save()
Write-Host "Applied all cell updates."